$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$vol = $ws.Range("A8").Characters(21, 2)
$vol.Text = "50"

$hdr = $ws.Range("C9").Characters(27, 9)
$hdr.Text = "12/8/2025"
$hdr2 = $ws.Range("C9").Characters(47, 9)
$hdr2.Text = "12/14/2025"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 66.666666666666
$ws.Range("C16").Value = 6
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 226
$ws.Range("J16").Value = 235
$ws.Range("K16").Value = -3.829787234042
$ws.Range("L16").Value = 49.668874172185
$ws.Range("M16").Value = 27.683615819209
$ws.Range("N16").Value = -82.120253164557
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = -7.692307692307
$ws.Range("I17").Value = 383
$ws.Range("J17").Value = 326
$ws.Range("K17").Value = 17.484662576687
$ws.Range("L17").Value = 56.326530612244
$ws.Range("M17").Value = 251.376146788991
$ws.Range("N17").Value = -8.373205741626
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -13.333333333333
$ws.Range("I18").Value = 221
$ws.Range("J18").Value = 229
$ws.Range("K18").Value = -3.493449781659
$ws.Range("L18").Value = 24.157303370786
$ws.Range("M18").Value = 93.859649122807
$ws.Range("N18").Value = -71.811224489795
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 694
$ws.Range("J19").Value = 759
$ws.Range("K19").Value = -8.563899868247
$ws.Range("L19").Value = 6.116207951070
$ws.Range("M19").Value = 42.213114754098
$ws.Range("N19").Value = -36.036866359447
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -50
$ws.Range("J20").Value = 94
$ws.Range("K20").Value = 7.446808510638
$ws.Range("L20").Value = 32.894736842105
$ws.Range("N20").Value = -85.893854748603
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 3.571428571428
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = -0.8
$ws.Range("I21").Value = 1643
$ws.Range("J21").Value = 1654
$ws.Range("K21").Value = -0.665054413542
$ws.Range("L21").Value = 25.419847328244
$ws.Range("M21").Value = 74.046610169491
$ws.Range("N21").Value = -61.790697674418
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 68
$ws.Range("J22").Value = 52
$ws.Range("K22").Value = 30.769230769230
$ws.Range("L22").Value = 25.925925925925
$ws.Range("M22").Value = -1.449275362318
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = -20
$ws.Range("M23").Value = 42.857142857142
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -24.242424242424
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 128
$ws.Range("H24").Value = -1.5625
$ws.Range("I24").Value = 1653
$ws.Range("J24").Value = 2104
$ws.Range("K24").Value = -21.43536121673
$ws.Range("L24").Value = -15.533980582524
$ws.Range("M24").Value = 13.296778615490
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -20
$ws.Range("F25").Value = 91
$ws.Range("G25").Value = 116
$ws.Range("H25").Value = -21.551724137931
$ws.Range("I25").Value = 1285
$ws.Range("J25").Value = 1870
$ws.Range("K25").Value = -31.283422459893
$ws.Range("L25").Value = -25.979262672811
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -40
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = -40
$ws.Range("I26").Value = 396
$ws.Range("J26").Value = 458
$ws.Range("K26").Value = -13.537117903930
$ws.Range("L26").Value = -1.246882793017
$ws.Range("M26").Value = 14.782608695652
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = 42.857142857142
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 20
$ws.Range("J28").Value = 56
$ws.Range("K28").Value = 33.928571428571
$ws.Range("L28").Value = 29.310344827586
$ws.Range("N29").Value = -86.842105263157
$ws.Range("N30").Value = -86.206896551724
$ws.Range("L31").Value = -57.142857142857

# --- Cells changing from text placeholder to numeric value (style 13 -> 14/15) ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = -100
$ws.Range("H15").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("G27").Value = 1
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("H27").Value = -100
$ws.Range("H27").NumberFormat = "#,##0.0;`"-`"#,##0.0"

# --- Cells changing from numeric value to text placeholder (style 14/15 -> 13) ---
# Use copy/paste-special (values then formats) from a stable reference cell so the
# existing shared-string + style index (13) is reused instead of allocating a new one.
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)